# This script reproduces the latest "cryptos list" price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, volume percentages, and any price strings
# that already look non-numeric) -- these can be assigned directly since Excel will
# keep them as text.
$ws.Range("D2").Value = '42.773.47'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '2.325.10'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("E6").Value = '  -3.87%  '
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -2.12%  '
$ws.Range("E10").Value = '  -4.58%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("E12").Value = '  -3.53%  '
$ws.Range("E13").Value = '  +1.87%  '
$ws.Range("E14").Value = '  -3.32%  '
$ws.Range("D15").Value = '2.693.13'
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '2.323.16'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").Value = '42.718.31'
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("E19").Value = '  -4.11%  '
$ws.Range("E20").Value = '  +1.10%  '
$ws.Range("D21").Value = '0.0₃0887'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("E24").Value = '  +0.98%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("E28").Value = '  +8.22%  '
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("E30").Value = '  -5.29%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  -16.21%  '
$ws.Range("E33").Value = '  -0.41%  '
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("E36").Value = '  -3.15%  '
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("E40").Value = '  +22.58%  '
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("D43").Value = '1.939.71'
$ws.Range("E43").Value = '  -2.99%  '
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("E45").Value = '  -5.25%  '
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").Value = '2.557.64'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  -1.64%  '
$ws.Range("E51").Value = '  +0.23%  '

# Some new "Price" values look exactly like plain decimal numbers (e.g. "302.03").
# Setting .Value directly on those would make Excel silently convert them to a
# numeric cell, which does not match the source data (it must stay text). To force
# text storage we temporarily enter a formula that evaluates to the literal text,
# then copy/paste-special as values only; that flattens the formula down to a plain
# text value while leaving the cells style untouched.
$riskyRefs = @("D4", "D5", "D6", "D10", "D11", "D12", "D14", "D17", "D19", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D34", "D36", "D44", "D45", "D46", "D50", "D51")
$riskyValues = @('0.999', '302.03', '94.08', '33.96', '0.0782', '18.73', '6.68', '0.789', '12.00', '6.14', '67.91', '235.62', '1.00', '2.42', '24.51', '2.23', '9.11', '31.45', '139.19', '17.58', '4.36', '0.0278', '10.17', '2.07', '52.84', '72.09')
for ($i = 0; $i -lt $riskyRefs.Length; $i++) {
    $cell = $ws.Range($riskyRefs[$i])
    $cell.Formula = '="' + $riskyValues[$i] + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}
$excel.CutCopyMode = $false

# Rebuild the calculation chain so no leftover references to the temporary formulas remain.
$excel.CalculateFullRebuild()